# Add files via upload
# Adds the new Time-Recording-Log entries (2019-11-16 .. 2019-12-09) that
# Tak Jaein ("탁재인") logged on sheet 6, rows 26-35, columns A:F.
# The TOTAL sheet's SUMIF formulas recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

# Make the new cells in column F look like the other "activity" text cells
# already on this sheet (Dotum font) instead of the default font the
# previously-blank cells had, by copying the format from F21 first.
$ws.Range("F21").Copy() | Out-Null
$ws.Range("F26:F35").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$rows = @(
    @{ Row=26; A=43785; B=0.875;               C=0.89583333333333337; D=0;   E=30;  F="자료조사" },
    @{ Row=27; A=43794; B=0.75;                C=0.875;               D=0;   E=180; F="개발" },
    @{ Row=28; A=43801; B=0.75;                C=0.875;               D=0;   E=180; F="개발" },
    @{ Row=29; A=43802; B=0.75;                C=0.95833333333333337; D=120; E=180; F="개발" },
    @{ Row=30; A=43804; B=0.83333333333333337; C=0.95833333333333337; D=30;  E=150; F="개발" },
    @{ Row=31; A=43805; B=0.54166666666666663; C=0.625;               D=0;   E=120; F="개발" },
    @{ Row=32; A=43807; B=0.33333333333333331; C=0.5;                 D=60;  E=180; F="개발" },
    @{ Row=33; A=43807; B=0.625;               C=0.83333333333333337; D=120; E=180; F="개발" },
    @{ Row=34; A=43807; B=0.875;               C=1;                   D=60;  E=120; F="개발" },
    @{ Row=35; A=43808; B=0.83333333333333337; C=1;                   D=60;  E=180; F="개발" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}

# Rows 34/35's date cells get a slightly different date number format
# (no ";@" text fallback) than the rest of the A column.
$ws.Range("A34:A35").NumberFormat = 'm"월"\ d"일"'

# The author's last action left the sheet scrolled down with F35 selected
# and this sheet as the active tab of the workbook.
$ws.Activate()
$ws.Range("F35").Select() | Out-Null
